# Revamp values of RQ2.1 and added entries from new papers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task-log rows appended after row 55 (row 56 left blank, matching the
# sheet's existing pattern of blank separator rows).
$ws.Range("A57").Value = "Note and update values of RQ2.1"
$ws.Range("B57").Value = 1
$ws.Range("C57").Value = 50

$ws.Range("A58").Value = "Updates the values of RQ2.1 and add values of new papers"
$ws.Range("B58").Value = 1
$ws.Range("C58").Value = 90

# Iterative calculation delta tweak reflected in workbook calcPr.
$excel.MaxChange = 0.0001

# Move/collapse the active selection the way the author's Excel session left it.
$ws.Range("A60").Select() | Out-Null
